$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.509.93"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "2.017.41"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.66"
$ws.Range("E5").Value = "  +6.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -1.74%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.08"
$ws.Range("E8").Value = "  -6.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0773"
$ws.Range("E10").Value = "  -3.32%  "

$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.41"
$ws.Range("E12").Value = "  -3.54%  "

$ws.Range("D13").Value = "2.315.33"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("E14").Value = "  -4.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.91"
$ws.Range("E15").Value = "  -7.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.25"
$ws.Range("E16").Value = "  -3.60%  "

$ws.Range("D17").Value = "2.022.92"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").Value = "37.453.10"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.65"
$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.84"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.68"
$ws.Range("E23").Value = "  +8.39%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.93"
$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.95"
$ws.Range("E27").Value = "  -4.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.65"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("E29").Value = "  -10.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("E30").Value = "  -1.44%  "

$ws.Range("E31").Value = "  -1.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0655"
$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("E33").Value = "  -3.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +1.66%  "

$ws.Range("E39").Value = "  -4.86%  "

$ws.Range("E40").Value = "  +4.69%  "

$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("E42").Value = "  -4.34%  "

$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("D44").Value = "1.393.78"
$ws.Range("E44").Value = "  +1.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.06"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.70"
$ws.Range("E46").Value = "  -5.76%  "

$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.04"
$ws.Range("E48").Value = "  -3.03%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D50").Value = "2.207.00"
$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -3.57%  "
